$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "collar_height" pressure-data column is being inserted ahead of the
# existing "Sample ID" block. Insert a column at L so the old L:AY range
# (Sample ID onward) shifts right to M:AZ, leaving J:L free for the new
# collar/submerged/exposed height trio.
$ws.Columns("L").Insert()

# Re-label the J:L header trio.
$ws.Range("J1").Value = "collar_height(cm)"
$ws.Range("K1").Value = "submerged_depth(cm)"
$ws.Range("L1").Value = "exposed_height(cm)"

# Start filling in the pressure-transducer readings for the existing rows.
$ws.Range("J2:J6").Value = 0
$ws.Range("K2:K6").Value = 0
$ws.Range("L2:L6").Value = 34.5

$ws.Range("L10").Select()
